$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new column D data: header "sconf_weight" plus numeric weights
$ws.Range("D1").Value = "sconf_weight"
$ws.Range("D2").Value = 0.7
$ws.Range("D3").Value = 1

# Update page setup (paper size / orientation)
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1

# Update the active selection to match the saved view state
$ws.Range("I9").Select()
